$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New calculation_time values (column F) written by the re-run of the metrics
# calculation pipeline the following day. Row -> new Excel serial datetime.
$newCalcTimes = @{
    2 = 45773.9126282363
    3 = 45773.91262875602
    4 = 45773.91262928458
    5 = 45773.9126297759
    6 = 45773.91263029215
    7 = 45773.9126308124
    8 = 45773.91263131925
    9 = 45773.91263186544
    10 = 45773.91263241286
    11 = 45773.91263292476
    12 = 45773.91263342176
    13 = 45773.9126339538
    14 = 45773.91263447842
    15 = 45773.91263496603
    16 = 45773.91263545559
    17 = 45773.91263594059
    18 = 45773.9126364302
    19 = 45773.91263692616
    20 = 45773.91263742017
    21 = 45773.9126379204
    22 = 45773.91263846097
    23 = 45773.91263906474
    24 = 45773.91263967872
    25 = 45773.91264050485
    26 = 45773.91264269219
    27 = 45773.91264364828
    28 = 45773.91264439347
    29 = 45773.9126453016
    30 = 45773.91264620119
    31 = 45773.91264696535
    32 = 45773.91264752953
    33 = 45773.91264802944
    34 = 45773.91264854317
    35 = 45773.91264905013
    36 = 45773.91264955245
    37 = 45773.91265006747
    38 = 45773.91265058192
    39 = 45773.91265107925
    40 = 45773.91265159703
    41 = 45773.91265211068
    42 = 45773.91265261287
    43 = 45773.91265311016
    44 = 45773.91265360954
    45 = 45773.91265410541
    46 = 45773.91265460493
    47 = 45773.91265509769
    48 = 45773.91265559488
    49 = 45773.91265615854
    50 = 45773.91265693743
    51 = 45773.91265758663
    52 = 45773.91265816594
    53 = 45773.91265886419
    54 = 45773.91265960685
    55 = 45773.91266012248
    56 = 45773.91266062044
    57 = 45773.9126611179
    58 = 45773.91266163401
    59 = 45773.91266163419
    60 = 45773.91266240962
    61 = 45773.91266325641
    62 = 45773.91266387749
    63 = 45773.91266387759
    64 = 45773.91266387761
    65 = 45773.91266387763
    66 = 45773.91266387764
    67 = 45773.91266387765
    68 = 45773.91266387767
    69 = 45773.91266387773
    70 = 45773.91266387775
    71 = 45773.91266387777
    72 = 45773.91266387777
    73 = 45773.91266387779
    74 = 45773.91266443131
    75 = 45773.91266443139
    76 = 45773.9126644314
    77 = 45773.91266443141
    78 = 45773.91266443142
    79 = 45773.91266443144
    80 = 45773.91266443145
    81 = 45773.91266443146
    82 = 45773.91266443147
    83 = 45773.91266443148
    84 = 45773.9126644315
    85 = 45773.9126644315
    86 = 45773.91266565309
    87 = 45773.9126662257
    88 = 45773.91266689355
    89 = 45773.91266691661
    90 = 45773.91266693729
    91 = 45773.91266695712
    92 = 45773.91266697707
    93 = 45773.91266699648
    94 = 45773.91266701606
    95 = 45773.91266703577
    96 = 45773.91266705455
    97 = 45773.91266707351
    98 = 45773.9126670922
    99 = 45773.91266711154
    100 = 45773.91266713069
    101 = 45773.91266714923
    102 = 45773.91266716787
    103 = 45773.91266718714
    104 = 45773.91266720785
    105 = 45773.91266722714
    106 = 45773.91266724555
    107 = 45773.91266726385
    108 = 45773.91266728217
    109 = 45773.91266730044
    110 = 45773.9126673189
    111 = 45773.91266733744
    112 = 45773.91266735609
    113 = 45773.91266737456
    114 = 45773.91266739279
    115 = 45773.91266741134
}

foreach ($row in $newCalcTimes.Keys) {
    $ws.Cells.Item($row, 6).Value = $newCalcTimes[$row]
}
